$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2094.5356
$ws.Range("J17").Value = 2094.5356
$ws.Range("L17").Value = 6283.6068
$ws.Range("N17").Value = -6619.6068
$ws.Range("H28").Value = 1120.125
$ws.Range("J28").Value = 829.1667
$ws.Range("L28").Value = 829.1667
$ws.Range("N28").Value = -1799.1667
$ws.Range("H55").Value = 71.85714
$ws.Range("I55").Value = 57.75
$ws.Range("J55").Value = 90.666664
$ws.Range("K55").Value = 57.75
$ws.Range("L55").Value = 90.666664
$ws.Range("M55").Value = 156.25
$ws.Range("N55").Value = -518.666664
$ws.Range("H62").Value = 6666.6665
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 6666.6665
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -56240
$ws.Range("H100").Value = 2908
$ws.Range("I100").Value = 2936
$ws.Range("J100").Value = 2880
$ws.Range("K100").Value = 2936
$ws.Range("L100").Value = 2880
$ws.Range("M100").Value = -2395
$ws.Range("N100").Value = -3962
$ws.Range("H127").Value = 1708.7858
$ws.Range("I127").Value = 648.5
$ws.Range("J127").Value = 1790.3462
$ws.Range("K127").Value = 1945.5
$ws.Range("L127").Value = 5371.0386
$ws.Range("M127").Value = 3014.5
$ws.Range("N127").Value = -15291.0386
$ws.Range("H137").Value = 1124901.9
$ws.Range("I137").Value = 2942294.8
$ws.Range("J137").Value = 1422.8
$ws.Range("K137").Value = 8826884.399999999
$ws.Range("L137").Value = 4268.4
$ws.Range("M137").Value = -8824334.399999999
$ws.Range("N137").Value = -9368.4
$ws.Range("H138").Value = 2781529.5
$ws.Range("I138").Value = 2360.625
$ws.Range("J138").Value = 3209093.8
$ws.Range("K138").Value = 7081.875
$ws.Range("L138").Value = 9627281.399999999
$ws.Range("M138").Value = -1941.875
$ws.Range("N138").Value = -9637561.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5309485.5
$ws.Range("I32").Value = 5867992.5
$ws.Range("J32").Value = 3669
$ws.Range("K32").Value = 5867992.5
$ws.Range("L32").Value = 3669
$ws.Range("M32").Value = -5867705.5
$ws.Range("N32").Value = -4243
$ws.Range("H74").Value = 6632507.5
$ws.Range("I74").Value = 10040580
$ws.Range("J74").Value = 78522.30499999999
$ws.Range("K74").Value = 10040580
$ws.Range("L74").Value = 78522.30499999999
$ws.Range("M74").Value = -10039706
$ws.Range("N74").Value = -80270.30499999999
$ws.Range("H77").Value = 6632507.5
$ws.Range("I77").Value = 10040580
$ws.Range("J77").Value = 78522.30499999999
$ws.Range("K77").Value = 50202900
$ws.Range("L77").Value = 392611.525
$ws.Range("M77").Value = -50198532
$ws.Range("N77").Value = -401347.525

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 26317840
$ws.Range("I105").Value = 31252072
$ws.Range("J105").Value = 1936.6666
$ws.Range("K105").Value = 31252072
$ws.Range("L105").Value = 1936.6666
$ws.Range("M105").Value = -31250325
$ws.Range("N105").Value = -5430.6666
$ws.Range("H107").Value = 3769.182
$ws.Range("I107").Value = 3552.2
$ws.Range("J107").Value = 3950
$ws.Range("K107").Value = 3552.2
$ws.Range("L107").Value = 3950
$ws.Range("M107").Value = -1632.2
$ws.Range("N107").Value = -7790
$ws.Range("H134").Value = 4093.2708
$ws.Range("I134").Value = 3232.3235
$ws.Range("J134").Value = 6184.143
$ws.Range("K134").Value = 9696.970499999999
$ws.Range("L134").Value = 18552.429
$ws.Range("M134").Value = -7161.970499999999
$ws.Range("N134").Value = -23622.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 720.0700000000001
$ws.Range("I31").Value = 470.09525
$ws.Range("J31").Value = 786.519
$ws.Range("K31").Value = 470.09525
$ws.Range("L31").Value = 786.519
$ws.Range("M31").Value = -175.09525
$ws.Range("N31").Value = -1376.519
$ws.Range("H34").Value = 720.0700000000001
$ws.Range("I34").Value = 470.09525
$ws.Range("J34").Value = 786.519
$ws.Range("K34").Value = 470.09525
$ws.Range("L34").Value = 786.519
$ws.Range("M34").Value = -268.09525
$ws.Range("N34").Value = -1190.519
$ws.Range("H62").Value = 4339
$ws.Range("I62").Value = 2231.6667
$ws.Range("J62").Value = 7500
$ws.Range("K62").Value = 2231.6667
$ws.Range("L62").Value = 7500
$ws.Range("M62").Value = -1607.6667
$ws.Range("N62").Value = -8748
$ws.Range("H65").Value = 4339
$ws.Range("I65").Value = 2231.6667
$ws.Range("J65").Value = 7500
$ws.Range("K65").Value = 11158.3335
$ws.Range("L65").Value = 37500
$ws.Range("M65").Value = -8038.333500000001
$ws.Range("N65").Value = -43740
$ws.Range("H94").Value = 3100.2
$ws.Range("I94").Value = 8066.6665
$ws.Range("J94").Value = 971.7143
$ws.Range("K94").Value = 8066.6665
$ws.Range("L94").Value = 971.7143
$ws.Range("M94").Value = -7615.6665
$ws.Range("N94").Value = -1873.7143
$ws.Range("H99").Value = 2690.476
$ws.Range("I99").Value = 2775
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 2775
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = -1277
$ws.Range("N99").Value = -3996
$ws.Range("H105").Value = 1236.4706
$ws.Range("I105").Value = 1268
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 1268
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 479
$ws.Range("N105").Value = -4494
$ws.Range("H107").Value = 598.2258
$ws.Range("I107").Value = 267.55554
$ws.Range("J107").Value = 1056.0769
$ws.Range("K107").Value = 267.55554
$ws.Range("L107").Value = 1056.0769
$ws.Range("M107").Value = 1652.44446
$ws.Range("N107").Value = -4896.0769
$ws.Range("H126").Value = 2690.476
$ws.Range("I126").Value = 2775
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 8325
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -5855
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 26826.45
$ws.Range("I132").Value = 1466.65
$ws.Range("J132").Value = 52186.25
$ws.Range("K132").Value = 4399.950000000001
$ws.Range("L132").Value = 156558.75
$ws.Range("M132").Value = -1869.950000000001
$ws.Range("N132").Value = -161618.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 611.1579
$ws.Range("I12").Value = 961.8182
$ws.Range("J12").Value = 129
$ws.Range("K12").Value = 2885.4546
$ws.Range("L12").Value = 387
$ws.Range("M12").Value = -2712.4546
$ws.Range("N12").Value = -733
$ws.Range("H58").Value = 1478.125
$ws.Range("I58").Value = 800
$ws.Range("J58").Value = 1500
$ws.Range("K58").Value = 2400
$ws.Range("L58").Value = 4500
$ws.Range("M58").Value = -2272
$ws.Range("N58").Value = -4756
$ws.Range("H113").Value = 526.3333
$ws.Range("I113").Value = 491.25
$ws.Range("J113").Value = 549.7222
$ws.Range("K113").Value = 1473.75
$ws.Range("L113").Value = 1649.1666
$ws.Range("M113").Value = 696.25
$ws.Range("N113").Value = -5989.1666
$ws.Range("H131").Value = 874.18866
$ws.Range("J131").Value = 945.1556
$ws.Range("L131").Value = 2835.4668
$ws.Range("N131").Value = -12915.4668
$ws.Range("H132").Value = 2113.4348
$ws.Range("I132").Value = 2228
$ws.Range("J132").Value = 2081.611
$ws.Range("K132").Value = 20052
$ws.Range("L132").Value = 18734.499
$ws.Range("M132").Value = -17522
$ws.Range("N132").Value = -23794.499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 606.6799999999999
$ws.Range("I107").Value = 432.6154
$ws.Range("J107").Value = 795.25
$ws.Range("K107").Value = 432.6154
$ws.Range("L107").Value = 795.25
$ws.Range("M107").Value = 1487.3846
$ws.Range("N107").Value = -4635.25
$ws.Range("H132").Value = 43862.688
$ws.Range("I132").Value = 37730.25
$ws.Range("J132").Value = 52448.1
$ws.Range("K132").Value = 113190.75
$ws.Range("L132").Value = 157344.3
$ws.Range("M132").Value = -110660.75
$ws.Range("N132").Value = -162404.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11911
$ws.Range("I40").Value = 20000
$ws.Range("J40").Value = 9888.75
$ws.Range("K40").Value = 20000
$ws.Range("L40").Value = 9888.75
$ws.Range("M40").Value = -19864
$ws.Range("N40").Value = -10160.75
$ws.Range("H46").Value = 1199
$ws.Range("I46").Value = 1199
$ws.Range("K46").Value = 1199
$ws.Range("M46").Value = -1011
$ws.Range("H122").Value = 3381.2222
$ws.Range("I122").Value = 3434.5625
$ws.Range("J122").Value = 2954.5
$ws.Range("K122").Value = 10303.6875
$ws.Range("L122").Value = 8863.5
$ws.Range("M122").Value = -7853.6875
$ws.Range("N122").Value = -13763.5
$ws.Range("H132").Value = 87394.086
$ws.Range("I132").Value = 3474
$ws.Range("J132").Value = 95023.17999999999
$ws.Range("K132").Value = 10422
$ws.Range("L132").Value = 285069.54
$ws.Range("M132").Value = -7892
$ws.Range("N132").Value = -290129.54
$ws.Range("H136").Value = 89509.664
$ws.Range("I136").Value = 37494.25
$ws.Range("J136").Value = 380796
$ws.Range("K136").Value = 112482.75
$ws.Range("L136").Value = 1142388
$ws.Range("M136").Value = -109932.75
$ws.Range("N136").Value = -1147488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3175.8333
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050
$ws.Range("H132").Value = 27967.334
$ws.Range("I132").Value = 18855.982
$ws.Range("J132").Value = 54821.844
$ws.Range("K132").Value = 56567.946
$ws.Range("L132").Value = 164465.532
$ws.Range("M132").Value = -54037.946
$ws.Range("N132").Value = -169525.532
